$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 5.1162812721297382
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.2017061611593038
$ws.Range("E2").ClearContents()

$ws.Range("B3").Value = 4.3431839514537645
$ws.Range("C3").Value = 5.3564084336275419
$ws.Range("D3").Value = 3.1758378627691446
$ws.Range("E3").Value = 8.0558901479131446

$ws.Range("B1:E3").Select()
